# Add python (Sqlite3 insert) and linux (load json from file) entries
# to the bottom of 工作表1, mirroring the existing "# Python" section rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$sqliteCode = @'
```python
#!/usr/bin/python
import sqlite3
conn = sqlite3.connect('test.db')
print "Opened database successfully";
conn.execute("INSERT INTO COMPANY (ID,NAME,AGE,ADDRESS,SALARY) \
      VALUES (1, 'Paul', 32, 'California', 20000.00 )");
conn.commit()
print "Records created successfully";
conn.close()
```
'@

$jsonCode = @'
import json
with open('strings.json') as json_data:
    d = json.load(json_data)
    print(d)
'@

$ws.Range("A14").Value = "# Python"
$ws.Range("B14").Value = "## Insert to Sqlite3"
$ws.Range("C14").Value = $sqliteCode

$ws.Range("A15").Value = "# Python"
$ws.Range("B15").Value = "## Load json from file"
$ws.Range("C15").Value = $jsonCode

# match the row/column styling used by the other detail rows in this sheet
# (column A: no-wrap style; columns B/C: wrap-text style) by copying the
# formatting straight from the row above, same as every other entry block
$ws.Range("A13").Copy()
$ws.Range("A14:A15").PasteSpecial(-4122)
$ws.Range("B13:C13").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)
$ws.Range("B15:C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Select() | Out-Null
